$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# Update the "quantity per variety" counters (columns K:O) for the rows
# whose coin varieties were newly acquired/recounted.

$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 1

$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 1

$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 1

$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 1

$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = 1
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 1

$ws.Range("K28").Value = 1
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 1
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 1

$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 2
$ws.Range("M29").Value = 1
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 1

$ws.Range("K30").Value = 1
$ws.Range("L30").Value = 1
$ws.Range("M30").Value = 1
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 1

$ws.Range("K31").Value = 1
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 1
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 1

$ws.Range("K32").Value = 1
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = 1
$ws.Range("N32").Value = 1
$ws.Range("O32").Value = 1

# Leave the selection on K27, matching the saved view's active cell.
$ws.Activate()
$ws.Range("K27").Select()
